$d = $word.ActiveDocument
$t = $d.Tables.Item(3)

function Replace-CellText($row, $col, $oldText, $newText) {
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    # Wrap:=0 (wdFindStop) and Replace:=1 (wdReplaceOne) keep this confined to the
    # cell's own Range instead of touching every matching run in the document.
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1)
    if (-not $found) {
        Write-Host "WARNING: replace not found at row $row col $col for '$oldText'"
    }
}

# Row 29 (Action header row of this block)
Replace-CellText 29 2 "Asset: Summary " "Asset: Fundamental Concept"
Replace-CellText 29 3 "Asset: Allocation" "Accounts: Income, Maintenance and Expenses "
Replace-CellText 29 4 "Asset: Vouchers, Bills, Taxes Check-List" "Investment and Business: Query Check-List"

# Row 30
Replace-CellText 30 2 "Asset: Fundamental Concept " "Asset: Inspection"
Replace-CellText 30 3 "Accounts: Income, Maintenance and Expenses " "Accounts: Unnecessary Expenses **  "

# Row 31 (col2 was empty, now gets new text)
$cell31_2 = $t.Cell(31, 2)
$rng31_2 = $cell31_2.Range
$rng31_2.End = $rng31_2.End - 1
$rng31_2.Text = "Asset: Summary and Valuation"

Replace-CellText 31 3 "Accounts: Unnecessary Expenses **  " "Asset: Vouchers, Bills, Taxes Check-List"

# Row 32
Replace-CellText 32 3 "Investment and Business: Query Check-List" "Asset: Allocation"

Write-Host "Done"
